$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename table column headers (renaming the header cell text also renames
# the corresponding ListObject/table column in-place, keeping the same
# column uid):
#  - Customer table column "Customer"  -> "CustomerKey"
#  - Sales table column    "Product"   -> "ProductKey"
#  - Product table column  "Product"   -> "ProductKey"
$ws.Range("E3").Value = "CustomerKey"
$ws.Range("F10").Value = "ProductKey"
$ws.Range("F19").Value = "ProductKey"

# Re-fit the columns whose header text changed so the displayed widths
# match the new header/content lengths.
$ws.Range("D1:F1").EntireColumn.AutoFit()
